# Tarea 4 - preguntas reflexivas: two small textual clarifications
# added within existing paragraphs of the "violencia domestica" answer.

$d = $word.ActiveDocument

# 1) Expand on what happens to the children ("congenitos") of an abusive
#    household: replace the single word with the fuller clarifying phrase.
$d.Content.Find.Execute(
    "congénitos",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "hijos (se tiene la creencia de que los hijos son propiedad de los padres, cayendo en violencia hacia los menores de edad)",
    2
)

# 2) Clarify that "ama de casa" is interpreted as a social norm.
$d.Content.Find.Execute(
    "solamente amas de casa, no cuentan",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "solamente amas de casa (interpretado como norma social), no cuentan",
    2
)
